$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @{
    "K"  = 2.14
    "M"  = 0.109
    "N"  = 0.002379912663755459
    "O"  = 0.05093457943925234
    "S"  = 0.109
    "U"  = 5.8
    "V"  = 0.1266375545851529
    "W"  = 0.3031161473087819
    "X"  = 0.07543308296429375
    "Y"  = 0.2276830643444882
    "AA" = 0.7734038193590691
    "AB" = 0.07543308296429375
    "AC" = 0.6979707363947754
    "AG" = -5.8
    "AJ" = -0.145
    "AK" = -1.611111111111111
    "AM" = -0.048
    "AP" = -2.735849056603773
    "AQ" = -41.04166666666666
}

foreach ($row in 2,3) {
    foreach ($col in $columns.Keys) {
        $ws.Range("$col$row").Value = $columns[$col]
    }
}
